# Generate Report for Archive
#
# 1) Status text update: every cell that shows "Ready for handoff" (the
#    Overview sheet's per-locale status columns, and the Status column on
#    the "zh-cn"/"de-de" detail sheets) now reads "In Translation".
# 2) Column width tightening: the two per-locale status columns on the
#    Overview sheet, and the Status column on each detail sheet, are
#    narrowed from ~17.22 characters to ~13.41 characters.

$wb = $excel.ActiveWorkbook

# --- 1) Replace the status text on every sheet -----------------------------
foreach ($ws in $wb.Worksheets) {
    $ws.Cells.Replace("Ready for handoff", "In Translation")
}

# --- 2) Narrow the "Status" columns ----------------------------------------
# Overview sheet: columns E (zh-cn) and F (de-de)
$overview = $wb.Worksheets.Item("Overview")
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

# zh-cn / de-de detail sheets: column C (Status)
$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = 12.5
